$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.25623881816864
$ws.Range("B1").Value = 2.474175214767456
$ws.Range("C1").Value = 4.935407161712646
$ws.Range("D1").Value = 3.007672786712646
$ws.Range("E1").Value = 1.112744450569153
